$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Summary block (rows 10-12): handle float/None-safe numeric inputs and
# recompute the marks summary.
# ---------------------------------------------------------------------------

# Row 10 ("No.") - give A10 the title style it was missing, refresh counts.
$ws.Range("A10").Style = "mtitleStyle"
$ws.Range("A10").HorizontalAlignment = -4108
$ws.Range("B10").Value = 24
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 3
$ws.Range("E10").Value = 28

# Row 11 ("Marking") - same title style fix; C11 becomes a real number
# instead of the inline string "-1" that used to break numeric consumers.
$ws.Range("A11").Style = "mtitleStyle"
$ws.Range("A11").HorizontalAlignment = -4108
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

# Row 12 ("Total")
$ws.Range("A12").Style = "mtitleStyle"
$ws.Range("A12").HorizontalAlignment = -4108
$ws.Range("B12").Value = 96
$ws.Range("C12").Value = -1
$ws.Range("E12").Value = "95/112"

# ---------------------------------------------------------------------------
# Third answer block (columns G/H) is no longer used - drop it completely.
# ---------------------------------------------------------------------------
$ws.Range("G15:H40").Clear()

# ---------------------------------------------------------------------------
# Second answer block (columns D/E): only rows 16-18 keep a "Student Ans"
# entry now, so clear everything below that and refresh the three that stay.
# ---------------------------------------------------------------------------
$ws.Range("D19:E40").Clear()

$ws.Range("D16").Style = "correctStyle"
$ws.Range("D16").HorizontalAlignment = -4108
$ws.Range("D16").Value = "Option A"

$ws.Range("D17").Style = "correctStyle"
$ws.Range("D17").HorizontalAlignment = -4108
$ws.Range("D17").Value = "Option C"

$ws.Range("D18").Style = "correctStyle"
$ws.Range("D18").HorizontalAlignment = -4108
$ws.Range("D18").Value = "Option D"

# ---------------------------------------------------------------------------
# First answer block (column A, "Student Ans"): fill in the student's
# answers. Matches the "Correct Ans" (column B) are shown with the green
# correctStyle; the one mismatch (row 36) is shown with the red
# incorrectStyle. Rows 24, 34 and 35 remain blank (not attempted).
# ---------------------------------------------------------------------------
$ws.Range("A16:A40").Style = "correctStyle"
$ws.Range("A16:A40").HorizontalAlignment = -4108

$ws.Range("A16").Value = "Option A"
$ws.Range("A17").Value = "Option D"
$ws.Range("A18").Value = "Option B"
$ws.Range("A19").Value = "Option C"
$ws.Range("A20").Value = "Option B"
$ws.Range("A21").Value = "Option C"
$ws.Range("A22").Value = "Option D"
$ws.Range("A23").Value = "Option D"
$ws.Range("A25").Value = "Option A"
$ws.Range("A26").Value = "Option C"
$ws.Range("A27").Value = "Option A"
$ws.Range("A28").Value = "Option D"
$ws.Range("A29").Value = "Option D"
$ws.Range("A30").Value = "Option B"
$ws.Range("A31").Value = "Option D"
$ws.Range("A32").Value = "Option C"
$ws.Range("A33").Value = "Option D"
$ws.Range("A37").Value = "Option A"
$ws.Range("A38").Value = "Option A"
$ws.Range("A39").Value = "Option D"
$ws.Range("A40").Value = "Option D"

# Rows 24, 34 and 35 were never answered - restore the blank/normal style.
$ws.Range("A24").Style = "normalStyle"
$ws.Range("A24").HorizontalAlignment = -4108

$ws.Range("A34").Style = "normalStyle"
$ws.Range("A34").HorizontalAlignment = -4108

$ws.Range("A35").Style = "normalStyle"
$ws.Range("A35").HorizontalAlignment = -4108

# Row 36 was answered wrong (correct answer is "Option A", student put
# "Option D") - flag it with the incorrect (red) style.
$ws.Range("A36").Style = "incorrectStyle"
$ws.Range("A36").HorizontalAlignment = -4108
$ws.Range("A36").Value = "Option D"
